# Applies the cell updates described by the commit diff to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = '67.960.39' },
    @{ Cell = "E2"; Value = '  -3.22%  ' },
    @{ Cell = "D3"; Value = '3.836.48' },
    @{ Cell = "E3"; Value = '  -2.61%  ' },
    @{ Cell = "E4"; Value = '  +0.01%  ' },
    @{ Cell = "D5"; Value = '599.98' },
    @{ Cell = "E5"; Value = '  -1.89%  ' },
    @{ Cell = "D6"; Value = '167.42' },
    @{ Cell = "E6"; Value = '  -2.54%  ' },
    @{ Cell = "D7"; Value = '3.836.31' },
    @{ Cell = "E7"; Value = '  -2.58%  ' },
    @{ Cell = "E8"; Value = '  +0.05%  ' },
    @{ Cell = "D9"; Value = '0.529' },
    @{ Cell = "E9"; Value = '  -2.03%  ' },
    @{ Cell = "D10"; Value = '0.164' },
    @{ Cell = "E10"; Value = '  -4.67%  ' },
    @{ Cell = "D11"; Value = '6.45' },
    @{ Cell = "E11"; Value = '  -0.07%  ' },
    @{ Cell = "D12"; Value = '0.458' },
    @{ Cell = "E12"; Value = '  -3.10%  ' },
    @{ Cell = "D13"; Value = '0.0000260' },
    @{ Cell = "E13"; Value = '  +0.45%  ' },
    @{ Cell = "D14"; Value = '36.94' },
    @{ Cell = "E14"; Value = '  -4.75%  ' },
    @{ Cell = "D15"; Value = '4.478.89' },
    @{ Cell = "E15"; Value = '  -2.71%  ' },
    @{ Cell = "D16"; Value = '3.830.49' },
    @{ Cell = "E16"; Value = '  -2.99%  ' },
    @{ Cell = "D17"; Value = '68.070.92' },
    @{ Cell = "E17"; Value = '  -3.15%  ' },
    @{ Cell = "D18"; Value = '18.18' },
    @{ Cell = "E18"; Value = '  -1.49%  ' },
    @{ Cell = "D19"; Value = '7.38' },
    @{ Cell = "E19"; Value = '  -4.02%  ' },
    @{ Cell = "E20"; Value = '  -0.99%  ' },
    @{ Cell = "D21"; Value = '10.99' },
    @{ Cell = "E21"; Value = '  -1.17%  ' },
    @{ Cell = "D22"; Value = '465.33' },
    @{ Cell = "E22"; Value = '  -6.70%  ' },
    @{ Cell = "D23"; Value = '0.732' },
    @{ Cell = "E23"; Value = '  -2.14%  ' },
    @{ Cell = "D24"; Value = '0.0000159' },
    @{ Cell = "E24"; Value = '  -4.74%  ' },
    @{ Cell = "D25"; Value = '82.88' },
    @{ Cell = "E25"; Value = '  -3.75%  ' },
    @{ Cell = "D26"; Value = '2.23' },
    @{ Cell = "E26"; Value = '  -3.47%  ' },
    @{ Cell = "D27"; Value = '12.10' },
    @{ Cell = "E27"; Value = '  -2.62%  ' },
    @{ Cell = "B28"; Value = 'RenderToken' },
    @{ Cell = "C28"; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' },
    @{ Cell = "D28"; Value = '10.05' },
    @{ Cell = "E28"; Value = '  -1.80%  ' },
    @{ Cell = "B29"; Value = 'Dai' },
    @{ Cell = "C29"; Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai' },
    @{ Cell = "D29"; Value = '0.998' },
    @{ Cell = "E29"; Value = '  -0.25%  ' },
    @{ Cell = "D30"; Value = '2.96' },
    @{ Cell = "E30"; Value = '  -2.06%  ' },
    @{ Cell = "D31"; Value = '3.984.46' },
    @{ Cell = "E31"; Value = '  -2.71%  ' },
    @{ Cell = "B32"; Value = 'NEARProtocol' },
    @{ Cell = "C32"; Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near' },
    @{ Cell = "D32"; Value = '7.61' },
    @{ Cell = "E32"; Value = '  -3.69%  ' },
    @{ Cell = "B33"; Value = 'ImmutableX' },
    @{ Cell = "C33"; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx' },
    @{ Cell = "D33"; Value = '2.31' },
    @{ Cell = "E33"; Value = '  -5.87%  ' },
    @{ Cell = "D34"; Value = '31.17' },
    @{ Cell = "E34"; Value = '  -3.87%  ' },
    @{ Cell = "D35"; Value = '9.57' },
    @{ Cell = "E35"; Value = '  -0.75%  ' },
    @{ Cell = "D36"; Value = '3.796.91' },
    @{ Cell = "E36"; Value = '  -2.77%  ' },
    @{ Cell = "E37"; Value = '  -4.11%  ' },
    @{ Cell = "D38"; Value = '3.58' },
    @{ Cell = "E38"; Value = '  +8.42%  ' },
    @{ Cell = "D39"; Value = '0.141' },
    @{ Cell = "E39"; Value = '  -0.99%  ' },
    @{ Cell = "E40"; Value = '  -2.73%  ' },
    @{ Cell = "D41"; Value = '5.91' },
    @{ Cell = "E41"; Value = '  -4.54%  ' },
    @{ Cell = "E42"; Value = '  -0.04%  ' },
    @{ Cell = "D43"; Value = '0.313' },
    @{ Cell = "E43"; Value = '  -5.15%  ' },
    @{ Cell = "D44"; Value = '1.98' },
    @{ Cell = "E44"; Value = '  -6.96%  ' },
    @{ Cell = "D45"; Value = '420.36' },
    @{ Cell = "E45"; Value = '  -4.64%  ' },
    @{ Cell = "B46"; Value = 'Cosmos' },
    @{ Cell = "C46"; Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom' },
    @{ Cell = "D46"; Value = '8.68' },
    @{ Cell = "E46"; Value = '  -0.20%  ' },
    @{ Cell = "B47"; Value = 'USDe' },
    @{ Cell = "C47"; Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde' },
    @{ Cell = "D47"; Value = '1.00' },
    @{ Cell = "E47"; Value = '  -0.03%  ' },
    @{ Cell = "B48"; Value = 'OKB' },
    @{ Cell = "C48"; Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb' },
    @{ Cell = "D48"; Value = '47.05' },
    @{ Cell = "E48"; Value = '  -2.49%  ' },
    @{ Cell = "B49"; Value = 'FLOKI' },
    @{ Cell = "C49"; Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki' },
    @{ Cell = "D49"; Value = '0.000289' },
    @{ Cell = "E49"; Value = '  +3.06%  ' },
    @{ Cell = "D50"; Value = '142.16' },
    @{ Cell = "E50"; Value = '  -0.86%  ' },
    @{ Cell = "E51"; Value = '  -3.54%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # Force text number-format first so numeric-looking strings (e.g. "1.00",
    # "0.0000260") are not silently coerced into doubles and lose their exact
    # textual representation (trailing zeros, thousands-style dot separators, etc).
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    # Restore the default/normal style so no stray formatting is introduced on
    # cells that originally had no explicit style (matches the source workbook).
    $cell.Style = "Normal"
}
